$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the manager/password pair in row 4 (October month data)
$ws.Range("A4").Value = "mngr353217"
$ws.Range("B4").Value = "ehadEru"

# Reflect the active selection left on cell B4 after the edit
$ws.Range("B4").Select()
